$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-code the ID (col A), Type (col B) and Site (col C) columns to reflect
# the new treatment/nursery/exploratory grouping scheme, and rename a few
# sites, ahead of building per-site graphs + summary stats.

$ws.Range("A2").Value = "E1"
$ws.Range("B2").Value = "exploratory"
$ws.Range("C2").Value = "Hell's Gate"

$ws.Range("A3").Value = "N2"
$ws.Range("B3").Value = "nursery"
$ws.Range("C3").Value = "Ten Pound Bay"

$ws.Range("A4").Value = "N1"
$ws.Range("B4").Value = "nursery"
$ws.Range("C4").Value = "York Island"

$ws.Range("A5").Value = "T3"
$ws.Range("B5").Value = "treatment"
$ws.Range("C5").Value = "Big Deep"

$ws.Range("A6").Value = "T4"
$ws.Range("B6").Value = "treatment"
$ws.Range("C6").Value = "Friar's Head Bay"

$ws.Range("A7").Value = "T5"
$ws.Range("B7").Value = "treatment"
$ws.Range("C7").Value = "Dansby's Beach"

$ws.Range("A8").Value = "T6"
$ws.Range("B8").Value = "treatment"
$ws.Range("C8").Value = "Club House"

$ws.Range("A9").Value = "E2"
$ws.Range("B9").Value = "exploratory"
$ws.Range("C9").Value = "Open Water"

$ws.Range("A10").Value = "T2"
$ws.Range("B10").Value = "treatment"
$ws.Range("C10").Value = "Green Isl. Anchorage"

$ws.Range("A11").Value = "T1"
$ws.Range("B11").Value = "treatment"
$ws.Range("C11").Value = "Nonsuch Anchorage"

# Move the active selection like the author left it.
$ws.Range("D7").Select()
